$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (autogluon) — previously empty inline strings, now populated.
$ws.Range("B3").Value = "0.399 (0.350 ± 0.021)"
$ws.Range("C3").Value = "00:02:39 (00:02:48 ± 00:00:07)"
$ws.Range("D3").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("E3").Value = "[]"

# F3 ("19") looks numeric; force it to stay text like its sibling cells
# (F4/F6/F8 etc. are all stored as text), then drop the quote-prefix style
# that forcing text via a leading apostrophe introduces so the cell keeps
# the sheet's default (unstyled) formatting.
$ws.Range("F3").Value = "'19"
$ws.Range("F3").Style = "Normal"

# Rows 4, 6, 8 — fix mojibake "Â±" -> "±" left over from a bad re-encode.
$ws.Range("B4").Value = "0.711 (0.677 ± 0.016)"
$ws.Range("C4").Value = "00:03:15 (00:03:49 ± 00:00:36)"
$ws.Range("D4").Value = "00:00:10 (00:00:10 ± 00:00:00)"

$ws.Range("B6").Value = "0.807 (0.775 ± 0.015)"
$ws.Range("C6").Value = "00:04:56 (00:05:01 ± 00:00:02)"
$ws.Range("D6").Value = "00:00:00 (00:00:02 ± 00:00:01)"

$ws.Range("B8").Value = "0.744 (0.689 ± 0.030)"
$ws.Range("C8").Value = "00:05:06 (00:09:55 ± 00:04:20)"
$ws.Range("D8").Value = "00:00:00 (00:00:00 ± 00:00:00)"
